$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "TC_Name" column (B), shifting
# TC_Name/Run/Status from B/C/D to C/D/E.
$ws.Columns("B").Insert()

# Header row
$ws.Range("B1").Value = "TC_path"

# New TC_path column values
$ws.Range("B2").Value = "./ui/tests/test_add_address.py"
$ws.Range("B3").Value = "./ui/tests/test_add_address.py"
$ws.Range("B4").Value = "./ui/tests/test_login.py"
$ws.Range("B5").Value = "./ui/tests/registration/test_registration.py"
$ws.Range("B6").Value = "./ui/tests/test_search.py"

# The swapped Status values for rows 3 & 4 (test_add_address -> no, test_login -> yes)
$ws.Range("D3").Value = "no"
$ws.Range("D4").Value = "yes"

# Fix formatting on the new column: the insert copied column A's style (which
# has a fill-applied variant on rows 5-6); rows 2-6 in the new TC_path column
# should match the plain style used elsewhere in the table (same as B2).
$ws.Range("B2").Copy()
$ws.Range("B5:B6").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Column widths: new TC_path column (B) wide, TC_Name column (now C) keeps
# its original best-fit width (already carried over automatically by Insert).
$ws.Columns("B").ColumnWidth = 38.8

# Update the selected cell to match the authored file.
$ws.Range("D4").Select() | Out-Null
